# Added suspension measurements file
# Appends a new "suspensionAngle" parameter row to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "suspensionAngle"
$ws.Range("B9").Value = 45
$ws.Range("C9").Value = "deg"

# Mirror the post-entry selection left behind after typing the new row
# (cursor drops to the next row, column B).
$ws.Range("B10").Select()
